$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# New log entries to append (rows 12 and 13)
$newRows = @(
    @("2025-09-10 18:49:37", "879928", "Trunkwala, Alaqmar Kutbuddin", "12"),
    @("2025-09-10 18:50:44", "879928", "Trunkwala, Alaqmar Kutbuddin", "12")
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    for ($col = 1; $col -le 4; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        # Force text storage so numeric-looking strings (IDs, grades) stay as text,
        # matching the rest of the sheet's data type.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col - 1]
        # Reset to the default style so no extra formatting/style index is introduced.
        $cell.Style = "Normal"
    }
}
